$d = $word.ActiveDocument
$ellipsis = [char]0x2026

# ------------------------------------------------------------------
# 1) "Josh <ellipsis>" -> "Josh" (drop the trailing space + ellipsis
#    that stood in for the missing last name).
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute("Josh " + $ellipsis, $false, $false, $false, $false, $false, `
                                  $true, 1, $false, "Josh", 2)
if (-not $found) {
    throw "Could not find the 'Josh <ellipsis>' placeholder text"
}

# ------------------------------------------------------------------
# 2) Locate the paragraph that now reads just "Josh" (after the
#    tabs) and append " Havens" as its own run, right after it -
#    matching the two separate <w:r> runs shown in the diff. A
#    single-character sentinel ("Z") is appended too; it gives us an
#    unambiguous anchor to drop the bookmark at in step 3 (placing a
#    bookmark exactly at a paragraph's trailing text boundary is
#    unreliable), and is removed again immediately afterwards.
# ------------------------------------------------------------------
$joshRange = $d.Content
$found2 = $joshRange.Find.Execute("Josh")
if (-not $found2) {
    throw "Could not find the 'Josh' run after stripping the ellipsis"
}
$para = $joshRange.Paragraphs(1)
$paraRange = $para.Range
$paraRange.InsertAfter(" HavensZ")

# ------------------------------------------------------------------
# 3) Re-seat the "_GoBack" bookmark: Word only ever keeps a single
#    "_GoBack", so adding it here removes it from wherever it used
#    to be (after "Communication/Scheduling:") and recreates it,
#    collapsed, immediately after " Havens" (i.e. right before the
#    "Z" sentinel we just added).
# ------------------------------------------------------------------
$havensRange = $d.Content
$foundHavens = $havensRange.Find.Execute("HavensZ")
if (-not $foundHavens) {
    throw "Could not find the 'HavensZ' sentinel text"
}
$sentinelStart = $havensRange.End - 1
$bookmarkSpot = $d.Range($sentinelStart, $sentinelStart)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null

# Remove the "Z" sentinel now that the bookmark anchors the correct
# position.
$sentinelRange = $d.Range($sentinelStart, $sentinelStart + 1)
$sentinelRange.Delete()
